$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "scikit-learn with GPU!"
$ws.Range("E12").Value = "https://tensorflow.blog/2022/03/29/scikit-learn-with-gpu/"

$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

$ws.Range("D32").Value = "내 코드를 테스트한다. (feat. pytest)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/358"

$ws.Range("D37").Value = "[Paper Review] BEIT: BERT Pre-Training of Image Transformers"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1964&mod=document&pageid=1"

$ws.Range("D46").Value = "어지럼증 (Dizziness) 감별진단"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/449"
